$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Binary Search 1")
$ws = $wb.Worksheets.Item("Binary Search 2")

# ---------------------------------------------------------------------------
# 1) Re-style row 9 (B9:F9) to the "highlighted" look already used by row 7
#    and by the matching row in the "Binary Search 1" sheet, by copying the
#    cell formatting (PasteSpecial -> formats only) from cells that already
#    carry the desired style.
# ---------------------------------------------------------------------------
$ws.Range("B7").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Range("C7").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("D9").PasteSpecial(-4122)

$ws1.Range("E10").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws1.Range("F10").Copy()
$ws.Range("F9").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) Add the new row 10 content
# ---------------------------------------------------------------------------
$ws.Range("B10").Value = 6
$ws.Range("C10").Value = "Binary S2 9"
$ws.Range("D10").Value = "Painter's partition problem"
$ws.Range("E10").Value = "Painter's Partition Problem - Problem | Scaler Academy"
$ws.Range("G10").Value = "Special modular division"

# Match formatting of row 10 to the equivalent "non highlighted" rows (row 8)
$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("D8").Copy()
$ws.Range("C10").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Rows.Item(10).RowHeight = 43.2

# ---------------------------------------------------------------------------
# 3) Hyperlink for the new E10 cell
# ---------------------------------------------------------------------------
$url = "https://www.scaler.com/academy/mentee-dashboard/class/30366/assignment/problems/271?navref=cl_tt_lst_nm"
$ws.Hyperlinks.Add($ws.Range("E10"), $url, [Type]::Missing, [Type]::Missing, $url)
$ws.Range("E10").Value = "Painter's Partition Problem - Problem | Scaler Academy"

# ---------------------------------------------------------------------------
# 4) Update view state to reflect the new selection
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("C10").Select()
